# Generate Report for Archive
# Moves the cd9fc48a-08cd-4161-afee-beef6704f90c row from the bottom of each
# status table up to just after ac4d99b3-...-0021f0923214 (i.e. row 5),
# pushing 10707ff1-... and b354880c-... down by one row, and flips the
# relocated row's Status from "Ready for handoff" to "In Translation" (its
# translation work is back in progress) on all three sheets: Overview,
# zh-cn, de-de.

$wb = $excel.ActiveWorkbook

function Set-RowValues($ws, $row, $values) {
    foreach ($col in $values.Keys) {
        $ws.Range($col + $row).Value = $values[$col]
    }
}

function Set-LinkDisplay($ws, $displays) {
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        foreach ($cellRef in $displays.Keys) {
            if ($addr -eq ('$' + $cellRef.Substring(0,1) + '$' + $cellRef.Substring(1))) {
                $h.TextToDisplay = $displays[$cellRef]
            }
        }
    }
}

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-RowValues $wsOverview 5 @{
    "A" = "cd9fc48a-08cd-4161-afee-beef6704f90c.md"
    "B" = "In Translation"
    "C" = "In Translation"
    "D" = "2016-03-24 22:42:03"
}
Set-RowValues $wsOverview 6 @{
    "A" = "10707ff1-eb6a-443c-9431-054b4f3aee71.md"
    "B" = "Ready for handoff"
    "C" = "Ready for handoff"
    "D" = "2016-03-24 22:38:54"
}
Set-RowValues $wsOverview 7 @{
    "A" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.md"
    "B" = "Ready for handoff"
    "C" = "Ready for handoff"
    "D" = "2016-03-24 22:42:45"
}

Set-LinkDisplay $wsOverview @{
    "A5" = "cd9fc48a-08cd-4161-afee-beef6704f90c.md"
    "A6" = "10707ff1-eb6a-443c-9431-054b4f3aee71.md"
    "A7" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.md"
}

# ---------- zh-cn sheet ----------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-RowValues $wsZhCn 5 @{
    "A" = "cd9fc48a-08cd-4161-afee-beef6704f90c.md"
    "C" = "In Translation"
    "D" = "cd9fc48a-08cd-4161-afee-beef6704f90c.3e289fc56ec5a62f889a45e3fe7009d72ac149fe.zh-cn.xlf"
    "E" = "2016-03-24 22:41:58"
}
Set-RowValues $wsZhCn 6 @{
    "A" = "10707ff1-eb6a-443c-9431-054b4f3aee71.md"
    "C" = "Ready for handoff"
    "D" = "10707ff1-eb6a-443c-9431-054b4f3aee71.11b37df63d6867dc439d82b97b26be9899d9454c.zh-cn.xlf"
    "E" = "2016-03-24 22:38:50"
}
Set-RowValues $wsZhCn 7 @{
    "A" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.md"
    "C" = "Ready for handoff"
    "D" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.c27ecce420d5d769275bc64fc40f150097be80d3.zh-cn.xlf"
    "E" = "2016-03-24 22:42:41"
}

Set-LinkDisplay $wsZhCn @{
    "A5" = "cd9fc48a-08cd-4161-afee-beef6704f90c.md"
    "D5" = "cd9fc48a-08cd-4161-afee-beef6704f90c.3e289fc56ec5a62f889a45e3fe7009d72ac149fe.zh-cn.xlf"
    "A6" = "10707ff1-eb6a-443c-9431-054b4f3aee71.md"
    "D6" = "10707ff1-eb6a-443c-9431-054b4f3aee71.11b37df63d6867dc439d82b97b26be9899d9454c.zh-cn.xlf"
    "A7" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.md"
    "D7" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.c27ecce420d5d769275bc64fc40f150097be80d3.zh-cn.xlf"
}

# ---------- de-de sheet ----------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-RowValues $wsDeDe 5 @{
    "A" = "cd9fc48a-08cd-4161-afee-beef6704f90c.md"
    "C" = "In Translation"
    "D" = "cd9fc48a-08cd-4161-afee-beef6704f90c.3e289fc56ec5a62f889a45e3fe7009d72ac149fe.de-de.xlf"
    "E" = "2016-03-24 22:42:03"
}
Set-RowValues $wsDeDe 6 @{
    "A" = "10707ff1-eb6a-443c-9431-054b4f3aee71.md"
    "C" = "Ready for handoff"
    "D" = "10707ff1-eb6a-443c-9431-054b4f3aee71.11b37df63d6867dc439d82b97b26be9899d9454c.de-de.xlf"
    "E" = "2016-03-24 22:38:54"
}
Set-RowValues $wsDeDe 7 @{
    "A" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.md"
    "C" = "Ready for handoff"
    "D" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.c27ecce420d5d769275bc64fc40f150097be80d3.de-de.xlf"
    "E" = "2016-03-24 22:42:45"
}

Set-LinkDisplay $wsDeDe @{
    "A5" = "cd9fc48a-08cd-4161-afee-beef6704f90c.md"
    "D5" = "cd9fc48a-08cd-4161-afee-beef6704f90c.3e289fc56ec5a62f889a45e3fe7009d72ac149fe.de-de.xlf"
    "A6" = "10707ff1-eb6a-443c-9431-054b4f3aee71.md"
    "D6" = "10707ff1-eb6a-443c-9431-054b4f3aee71.11b37df63d6867dc439d82b97b26be9899d9454c.de-de.xlf"
    "A7" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.md"
    "D7" = "b354880c-24dd-4b16-ac5e-af6f3853cae3.c27ecce420d5d769275bc64fc40f150097be80d3.de-de.xlf"
}

Write-Host "Report regenerated: cd9fc48a moved to row 5 on all sheets."
